$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Good Morning" text in E8 with "GIT UPDATE" (this also drops
# the now-unused "Good Morning" shared string and appends the new one).
$ws.Range("E8").Value = "GIT UPDATE"

# Update the sheet's active selection to E8, matching the saved view state.
$ws.Range("E8").Select()
